$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 77, shifting existing rows 77-100 down to 78-101
$ws.Rows.Item(77).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Fill in the new row 77 with the latest weekly price entry
$ws.Cells.Item(77, 1).Value = 5
$ws.Cells.Item(77, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(77, 3).Value = "Maule"
$ws.Cells.Item(77, 4).Value = 45119
$ws.Cells.Item(77, 5).Value = 7
$ws.Cells.Item(77, 6).Value = "Fruta"
$ws.Cells.Item(77, 7).Value = 100104
$ws.Cells.Item(77, 8).Value = "Frutos de pepita"
$ws.Cells.Item(77, 9).Value = 100104003
$ws.Cells.Item(77, 10).Value = "Membrillo"
$ws.Cells.Item(77, 11).Value = "Champion"
$ws.Cells.Item(77, 12).Value = "Especial"
$ws.Cells.Item(77, 13).Value = 210
$ws.Cells.Item(77, 14).Value = 12000
$ws.Cells.Item(77, 15).Value = 12000
$ws.Cells.Item(77, 16).Value = 12000
$ws.Cells.Item(77, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(77, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(77, 19).Value = 667
$ws.Cells.Item(77, 20).Value = 18
